# 18:46:43.27 / 2025/01/19 周日
# Add a new "PhaseMask" text box to slide 1, as a new top-level shape
# (sibling of the big existing group), matching the authored OOXML:
#   <p:sp> id=2 name="文本框 1", rect 6029960,1186815 1790700x460375 EMU,
#   no fill, hairline-but-invisible centered single line, centered text
#   "PhaseMask" at 24pt with shape autofit / square wrap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points (914400 EMU per inch, 12700 EMU per point)
$left   = 6029960 / 12700
$top    = 1186815 / 12700
$width  = 1790700 / 12700
$height = 460375  / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "文本框 1"

$tr = $shp.TextFrame.TextRange
$tr.Text = "PhaseMask"
$tr.Font.Size = 24
$tr.ParagraphFormat.Alignment = 2

$shp.Fill.Visible = 0

$shp.Line.Visible = 0
$shp.Line.Weight = 2.25
$shp.Line.DashStyle = 1
$shp.Line.Style = 1

$shp.TextFrame.WordWrap = -1
$shp.TextFrame.AutoSize = 1

# AutoSize recomputes Height from font metrics (a few EMU off the
# canonical authored value) - put back the exact authored height.
# (Left/Top/Width are untouched by AutoSize and already exact from
# AddTextbox; re-assigning Left/Top goes through a lossy point
# round-trip in this host and must be avoided.)
$shp.Height = $height

Write-Output ("Added shape '" + $shp.Name + "' ; slide now has " + $s.Shapes.Count + " shapes")
